# "Generate Report for Archive"
#
# Refresh the localization-status report: the zh-cn / de-de hand-off status
# moves on from "Ready for handoff" to "In Translation" everywhere it is
# shown (the Overview summary columns, and the Status column on each
# language sheet), and the now-shorter status text lets those columns
# re-autofit a bit narrower.

$wb = $excel.ActiveWorkbook

$newStatus  = "In Translation"
$newWidth   = 12.5   # ColumnWidth (characters) -> narrower autofit for the shorter text

# Overview sheet: zh-cn (E) and de-de (F) status columns
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

# zh-cn sheet: Status column (C)
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

# de-de sheet: Status column (C)
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Columns.Item(3).ColumnWidth = $newWidth
